# "frissen felvett mentés eredményekbe emelése"
# A newly-added entry is lifted into the results sheet: append a new row
# (row 8) to Sheet1 with the next Verseny_ID, mirroring the existing
# (empty) Verseny_start / Verseny_end / Szervezo columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 8

$ws.Range("A$newRow").Value = "VID_00007"

# The other three columns stay empty for every existing row, but they are
# still present as empty *text* cells (not blank/number cells). A leading
# apostrophe forces an empty-text entry instead of clearing the cell, and
# resetting the Style afterwards drops the "quote prefix" formatting that
# Excel would otherwise remember for that cell.
$ws.Range("B$newRow").Value = "'"
$ws.Range("C$newRow").Value = "'"
$ws.Range("D$newRow").Value = "'"
$ws.Range("B$newRow`:D$newRow").Style = "Normal"
